# "save data done + era data updated"
# Add a new "Save" column (H) to the sheet:
#   - H1: header label "Save", styled like the other header cells (bold,
#         centered, bordered) by copying the formatting from the existing
#         G1 header cell.
#   - H2: the new data point's Save value (0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy G1's formatting (bold/centered/bordered header style) onto H1 so the
# new header cell matches the rest of row 1 without introducing a new style.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Set the new header text and data value.
$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0
